$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) pitstop sheet: add tyre_before / tyre_after columns (H, I)
# ---------------------------------------------------------------------------
$wsPit = $wb.Worksheets.Item("pitstop")

$wsPit.Range("H1").Value = "tyre_before"
$wsPit.Range("I1").Value = "tyre_after"

$tyreBefore = @(4,4,3,4,4,3,3,3,4,4,4,3,4,4,4,3,4,4,4,3,3,4,3,3,3,4)
$tyreAfter  = @(3,3,4,3,3,4,4,4,4,3,3,4,3,3,3,4,3,3,3,3,4,3,3,3,4,3)

for ($i = 0; $i -lt $tyreBefore.Length; $i++) {
    $row = $i + 2
    $wsPit.Cells.Item($row, 8).Value = $tyreBefore[$i]
    $wsPit.Cells.Item($row, 9).Value = $tyreAfter[$i]
}

[void]$wsPit.Range("H28").Select()

# ---------------------------------------------------------------------------
# 2) add "weather" worksheet at the end
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsWeather = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsWeather.Name = "weather"

$wsWeather.Range("A1").Value = "Skycondition"
$wsWeather.Range("B1").Value = "Clear"
$wsWeather.Range("A2").Value = "Temperature"
$wsWeather.Range("B2").Value = "69.98°F"
$wsWeather.Range("A3").Value = "Humidity"
$wsWeather.Range("B3").Value = 0.44
$wsWeather.Range("B3").NumberFormat = "0%"
$wsWeather.Range("A4").Value = "Wind speed"
$wsWeather.Range("B4").Value = "18.12 mph"
$wsWeather.Range("A5").Value = "Wind bearing"
$wsWeather.Range("B5").Value = "171°"

[void]$wsWeather.Range("D4").Select()

# ---------------------------------------------------------------------------
# 3) add "altitude" worksheet at the end
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsAltitude = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsAltitude.Name = "altitude"

$wsAltitude.Range("A1").Value = "delta"
$wsAltitude.Range("B1").Value = 26.8

[void]$wsAltitude.Range("B2").Select()
